# Update cryptos list values (price "D" column and volume/1h "E" column)
# to match the latest GitHub Actions scrape. Row 41/42 additionally swap
# the Coin name/Link (B/C) because Kaspa and EnergySwap traded ranking
# positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '96.417.06'
$ws.Cells.Item(2, 5).Value = '  -1.00%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.713.09'
$ws.Cells.Item(3, 5).Value = '  +3.33%  '

$ws.Cells.Item(4, 5).Value = '  +0.06%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '238.12'
$ws.Cells.Item(5, 5).Value = '  -2.19%  '

$ws.Cells.Item(6, 5).Value = '  +6.94%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '655.19'
$ws.Cells.Item(7, 5).Value = '  +0.39%  '

$ws.Cells.Item(8, 5).Value = '  -0.76%  '

$ws.Cells.Item(9, 5).Value = '  +1.23%  '

$ws.Cells.Item(10, 5).Value = '  +0.05%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '3.708.16'
$ws.Cells.Item(11, 5).Value = '  +3.27%  '

$ws.Cells.Item(12, 5).Value = '  +0.00%  '

$ws.Cells.Item(13, 5).Value = '  +0.39%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.84'
$ws.Cells.Item(14, 5).Value = '  +5.81%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '4.406.78'
$ws.Cells.Item(15, 5).Value = '  +3.37%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.0000267'
$ws.Cells.Item(16, 5).Value = '  +2.75%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '96.241.81'
$ws.Cells.Item(17, 5).Value = '  -0.86%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '8.98'
$ws.Cells.Item(18, 5).Value = '  +15.85%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '3.708.05'
$ws.Cells.Item(19, 5).Value = '  +3.33%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '19.10'
$ws.Cells.Item(20, 5).Value = '  +4.42%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '12.76'
$ws.Cells.Item(21, 5).Value = '  +1.05%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.525'
$ws.Cells.Item(22, 5).Value = '  -1.31%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '523.59'
$ws.Cells.Item(23, 5).Value = '  +1.11%  '

$ws.Cells.Item(24, 5).Value = '  -0.30%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '7.02'
$ws.Cells.Item(25, 5).Value = '  +0.66%  '

$ws.Cells.Item(26, 5).Value = '  -1.55%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '102.56'
$ws.Cells.Item(27, 5).Value = '  -0.60%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '13.44'
$ws.Cells.Item(28, 5).Value = '  +0.85%  '

$ws.Cells.Item(29, 5).Value = '  -7.74%  '

$ws.Cells.Item(30, 5).Value = '  +3.36%  '

$ws.Cells.Item(31, 5).Value = '  +2.39%  '

$ws.Cells.Item(32, 5).Value = '  +0.25%  '

$ws.Cells.Item(33, 5).Value = '  +10.29%  '

$ws.Cells.Item(34, 5).Value = '  -2.57%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '669.83'
$ws.Cells.Item(35, 5).Value = '  +8.47%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '32.77'
$ws.Cells.Item(36, 5).Value = '  +2.89%  '

$ws.Cells.Item(37, 5).Value = '  +0.32%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.598'
$ws.Cells.Item(38, 5).Value = '  +2.18%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '8.86'
$ws.Cells.Item(39, 5).Value = '  +0.68%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '7.10'
$ws.Cells.Item(40, 5).Value = '  +16.13%  '

$ws.Cells.Item(41, 2).Value = 'Kaspa'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.161'
$ws.Cells.Item(41, 5).Value = '  +4.33%  '

$ws.Cells.Item(42, 2).Value = 'EnergySwap'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '40.51'
$ws.Cells.Item(42, 5).Value = '  +23.49%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.977'
$ws.Cells.Item(43, 5).Value = '  +5.11%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.97'
$ws.Cells.Item(44, 5).Value = '  +2.25%  '

$ws.Cells.Item(45, 5).Value = '  +0.04%  '

$ws.Cells.Item(46, 5).Value = '  +1.50%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.436'
$ws.Cells.Item(47, 5).Value = '  -2.65%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.31'
$ws.Cells.Item(48, 5).Value = '  -1.39%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '23.60'
$ws.Cells.Item(49, 5).Value = '  -0.24%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '8.60'
$ws.Cells.Item(50, 5).Value = '  -1.76%  '

$ws.Cells.Item(51, 5).Value = '  +2.38%  '
